$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I3").Value = -0.1805290139085289
$ws.Range("J3").Value = 0.2074456135354681
$ws.Range("K3").Value = -0.7395785195415453
$ws.Range("L3").Value = 3.052775502476929
